$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shift D:K -> F:M)
$ws.Columns("D:E").Insert()

# Copy number formatting from the (now shifted) old D:K range (F:M) into new D:E
$ws.Range("F7:F102").Copy() | Out-Null
$ws.Range("D7:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Write quarterly financial data, columns D..M = 10 quarters
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43281
$ws.Range("G7").Value = 43190
$ws.Range("H7").Value = 43100
$ws.Range("I7").Value = 43008
$ws.Range("J7").Value = 42916
$ws.Range("K7").Value = 42825
$ws.Range("L7").Value = 42735
$ws.Range("M7").Value = 42643
$ws.Range("D8").Value = 243200
$ws.Range("E8").Value = 278000
$ws.Range("F8").Value = 243700
$ws.Range("G8").Value = 197200
$ws.Range("H8").Value = 201800
$ws.Range("I8").Value = 145900
$ws.Range("J8").Value = 168000
$ws.Range("K8").Value = 175000
$ws.Range("L8").Value = 167000
$ws.Range("M8").Value = 144700
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = "NA"
$ws.Range("H9").Value = "NA"
$ws.Range("I9").Value = "NA"
$ws.Range("J9").Value = "NA"
$ws.Range("K9").Value = "NA"
$ws.Range("L9").Value = "NA"
$ws.Range("M9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "NA"
$ws.Range("G10").Value = "NA"
$ws.Range("H10").Value = "NA"
$ws.Range("I10").Value = "NA"
$ws.Range("J10").Value = "NA"
$ws.Range("K10").Value = "NA"
$ws.Range("L10").Value = "NA"
$ws.Range("M10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("L12").Value = "NA"
$ws.Range("M12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "NA"
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = -225900
$ws.Range("M14").Value = 46800
$ws.Range("D15").Value = 63700
$ws.Range("E15").Value = 60700
$ws.Range("F15").Value = 54500
$ws.Range("G15").Value = 47700
$ws.Range("H15").Value = 48900
$ws.Range("I15").Value = 44500
$ws.Range("J15").Value = 48200
$ws.Range("K15").Value = 46500
$ws.Range("L15").Value = 49100
$ws.Range("M15").Value = 70400
$ws.Range("D17").Value = 164600
$ws.Range("E17").Value = 161100
$ws.Range("F17").Value = 150700
$ws.Range("G17").Value = 136100
$ws.Range("H17").Value = 134900
$ws.Range("I17").Value = 120700
$ws.Range("J17").Value = 125300
$ws.Range("K17").Value = 133600
$ws.Range("L17").Value = -81900
$ws.Range("M17").Value = 203400
$ws.Range("D18").Value = 78600
$ws.Range("E18").Value = 116900
$ws.Range("F18").Value = 92900
$ws.Range("G18").Value = 61100
$ws.Range("H18").Value = 66800
$ws.Range("I18").Value = 25200
$ws.Range("J18").Value = 42700
$ws.Range("K18").Value = 41400
$ws.Range("L18").Value = 248900
$ws.Range("M18").Value = -58800
$ws.Range("D20").Value = 168100
$ws.Range("E20").Value = -34600
$ws.Range("F20").Value = -74300
$ws.Range("G20").Value = -22900
$ws.Range("H20").Value = -33000
$ws.Range("I20").Value = -12400
$ws.Range("J20").Value = 91200
$ws.Range("K20").Value = 47200
$ws.Range("L20").Value = -32600
$ws.Range("M20").Value = -500
$ws.Range("D21").Value = 310300
$ws.Range("E21").Value = 143000
$ws.Range("F21").Value = 73000
$ws.Range("G21").Value = 85800
$ws.Range("H21").Value = 82600
$ws.Range("I21").Value = 57300
$ws.Range("J21").Value = 182100
$ws.Range("K21").Value = 135100
$ws.Range("L21").Value = 265400
$ws.Range("M21").Value = 11100
$ws.Range("D22").Value = 7300
$ws.Range("E22").Value = 6400
$ws.Range("F22").Value = 6900
$ws.Range("G22").Value = 6800
$ws.Range("H22").Value = 7200
$ws.Range("I22").Value = 6400
$ws.Range("J22").Value = 7600
$ws.Range("K22").Value = 7800
$ws.Range("L22").Value = 8600
$ws.Range("M22").Value = 100
$ws.Range("D23").Value = 239300
$ws.Range("E23").Value = 75900
$ws.Range("F23").Value = 11700
$ws.Range("G23").Value = 31400
$ws.Range("H23").Value = 26600
$ws.Range("I23").Value = 6300
$ws.Range("J23").Value = 126300
$ws.Range("K23").Value = 80800
$ws.Range("L23").Value = 207700
$ws.Range("M23").Value = -59400
$ws.Range("D24").Value = 53800
$ws.Range("E24").Value = 11300
$ws.Range("F24").Value = 2500
$ws.Range("G24").Value = 9300
$ws.Range("H24").Value = 15200
$ws.Range("I24").Value = -5700
$ws.Range("J24").Value = 30000
$ws.Range("K24").Value = 22200
$ws.Range("L24").Value = -437800
$ws.Range("M24").Value = 17900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("D26").Value = 185500
$ws.Range("E26").Value = 64700
$ws.Range("F26").Value = 9200
$ws.Range("G26").Value = 22100
$ws.Range("H26").Value = 11400
$ws.Range("I26").Value = 12000
$ws.Range("J26").Value = 96200
$ws.Range("K26").Value = 58600
$ws.Range("L26").Value = 645600
$ws.Range("M26").Value = -77400
$ws.Range("D27").Value = 185500
$ws.Range("E27").Value = 64700
$ws.Range("F27").Value = 9200
$ws.Range("G27").Value = 22100
$ws.Range("H27").Value = 11400
$ws.Range("I27").Value = 12000
$ws.Range("J27").Value = 96200
$ws.Range("K27").Value = 58600
$ws.Range("L27").Value = 645600
$ws.Range("M27").Value = -77400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("D32").Value = -168100
$ws.Range("E32").Value = 34600
$ws.Range("F32").Value = 74300
$ws.Range("G32").Value = 22900
$ws.Range("H32").Value = 33000
$ws.Range("I32").Value = 12400
$ws.Range("J32").Value = -91200
$ws.Range("K32").Value = -47200
$ws.Range("L32").Value = 32600
$ws.Range("M32").Value = 500
$ws.Range("D33").Value = 185500
$ws.Range("E33").Value = 64700
$ws.Range("F33").Value = 9200
$ws.Range("G33").Value = 22100
$ws.Range("H33").Value = 11400
$ws.Range("I33").Value = 12000
$ws.Range("J33").Value = 96200
$ws.Range("K33").Value = 58600
$ws.Range("L33").Value = 645600
$ws.Range("M33").Value = -77400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("D35").Value = 185500
$ws.Range("E35").Value = 64700
$ws.Range("F35").Value = 9200
$ws.Range("G35").Value = 22100
$ws.Range("H35").Value = 11400
$ws.Range("I35").Value = 12000
$ws.Range("J35").Value = 96200
$ws.Range("K35").Value = 58600
$ws.Range("L35").Value = 645600
$ws.Range("M35").Value = -77400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43281
$ws.Range("G38").Value = 43190
$ws.Range("H38").Value = 43100
$ws.Range("I38").Value = 43008
$ws.Range("J38").Value = 42916
$ws.Range("K38").Value = 42825
$ws.Range("L38").Value = 42735
$ws.Range("M38").Value = 42643
$ws.Range("D41").Value = 270400
$ws.Range("E41").Value = 258700
$ws.Range("F41").Value = 268200
$ws.Range("G41").Value = 295000
$ws.Range("H41").Value = 257900
$ws.Range("I41").Value = 259800
$ws.Range("J41").Value = 286600
$ws.Range("K41").Value = 3200
$ws.Range("L41").Value = 1000
$ws.Range("M41").Value = 57600
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("D43").Value = 149100
$ws.Range("E43").Value = 187700
$ws.Range("F43").Value = 172900
$ws.Range("G43").Value = 101200
$ws.Range("H43").Value = 97200
$ws.Range("I43").Value = 64200
$ws.Range("J43").Value = 70900
$ws.Range("K43").Value = 71600
$ws.Range("L43").Value = 88600
$ws.Range("M43").Value = 63100
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("D45").Value = 50700
$ws.Range("E45").Value = 3200
$ws.Range("F45").Value = 2200
$ws.Range("G45").Value = 5100
$ws.Range("H45").Value = 7300
$ws.Range("I45").Value = 16300
$ws.Range("J45").Value = 29900
$ws.Range("K45").Value = 313800
$ws.Range("L45").Value = 306300
$ws.Range("M45").Value = 13100
$ws.Range("D46").Value = 470300
$ws.Range("E46").Value = 449600
$ws.Range("F46").Value = 443400
$ws.Range("G46").Value = 401300
$ws.Range("H46").Value = 362300
$ws.Range("I46").Value = 340200
$ws.Range("J46").Value = 387400
$ws.Range("K46").Value = 388600
$ws.Range("L46").Value = 395900
$ws.Range("M46").Value = 133800
$ws.Range("D47").Value = 20200
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 38200
$ws.Range("H47").Value = 37300
$ws.Range("I47").Value = "NA"
$ws.Range("J47").Value = "NA"
$ws.Range("K47").Value = "NA"
$ws.Range("L47").Value = "NA"
$ws.Range("M47").Value = "NA"
$ws.Range("D48").Value = 972700
$ws.Range("E48").Value = 927700
$ws.Range("F48").Value = 853000
$ws.Range("G48").Value = 756000
$ws.Range("H48").Value = 669800
$ws.Range("I48").Value = 609900
$ws.Range("J48").Value = 579400
$ws.Range("K48").Value = 595700
$ws.Range("L48").Value = 567300
$ws.Range("M48").Value = 545400
$ws.Range("D49").Value = 487300
$ws.Range("E49").Value = 479200
$ws.Range("F49").Value = 481700
$ws.Range("G49").Value = 478900
$ws.Range("H49").Value = 475500
$ws.Range("I49").Value = 474400
$ws.Range("J49").Value = 480000
$ws.Range("K49").Value = 499400
$ws.Range("L49").Value = 500600
$ws.Range("M49").Value = 497100
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("D52").Value = 370100
$ws.Range("E52").Value = 408700
$ws.Range("F52").Value = 423900
$ws.Range("G52").Value = 421500
$ws.Range("H52").Value = 424200
$ws.Range("I52").Value = 475000
$ws.Range("J52").Value = 491200
$ws.Range("K52").Value = 544300
$ws.Range("L52").Value = 563400
$ws.Range("M52").Value = 126000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("D54").Value = 2320700
$ws.Range("E54").Value = 2265100
$ws.Range("F54").Value = 2202000
$ws.Range("G54").Value = 2096000
$ws.Range("H54").Value = 1969100
$ws.Range("I54").Value = 1899500
$ws.Range("J54").Value = 1938000
$ws.Range("K54").Value = 2028000
$ws.Range("L54").Value = 2027200
$ws.Range("M54").Value = 1302300
$ws.Range("D57").Value = 130000
$ws.Range("E57").Value = 124700
$ws.Range("F57").Value = 122100
$ws.Range("G57").Value = 97100
$ws.Range("H57").Value = 87200
$ws.Range("I57").Value = 64200
$ws.Range("J57").Value = 53100
$ws.Range("K57").Value = 65000
$ws.Range("L57").Value = 61200
$ws.Range("M57").Value = 55900
$ws.Range("D58").Value = 44700
$ws.Range("E58").Value = 43500
$ws.Range("F58").Value = 43800
$ws.Range("G58").Value = 21100
$ws.Range("H58").Value = 20600
$ws.Range("I58").Value = 20400
$ws.Range("J58").Value = 21200
$ws.Range("K58").Value = 22500
$ws.Range("L58").Value = 22700
$ws.Range("M58").Value = 22200
$ws.Range("D59").Value = 89100
$ws.Range("E59").Value = 207500
$ws.Range("F59").Value = 186900
$ws.Range("G59").Value = 144500
$ws.Range("H59").Value = 95100
$ws.Range("I59").Value = 82400
$ws.Range("J59").Value = 83700
$ws.Range("K59").Value = 97700
$ws.Range("L59").Value = 104400
$ws.Range("M59").Value = 78600
$ws.Range("D60").Value = 263700
$ws.Range("E60").Value = 375700
$ws.Range("F60").Value = 352800
$ws.Range("G60").Value = 262700
$ws.Range("H60").Value = 202900
$ws.Range("I60").Value = 166900
$ws.Range("J60").Value = 158100
$ws.Range("K60").Value = 185200
$ws.Range("L60").Value = 188300
$ws.Range("M60").Value = 156700
$ws.Range("D61").Value = 473900
$ws.Range("E61").Value = 448600
$ws.Range("F61").Value = 456400
$ws.Range("G61").Value = 491200
$ws.Range("H61").Value = 479800
$ws.Range("I61").Value = 476200
$ws.Range("J61").Value = 494600
$ws.Range("K61").Value = 549000
$ws.Range("L61").Value = 567900
$ws.Range("M61").Value = 537900
$ws.Range("D62").Value = 93900
$ws.Range("E62").Value = 129500
$ws.Range("F62").Value = 121000
$ws.Range("G62").Value = 100700
$ws.Range("H62").Value = 95000
$ws.Range("I62").Value = 82500
$ws.Range("J62").Value = 82400
$ws.Range("K62").Value = 120300
$ws.Range("L62").Value = 149000
$ws.Range("M62").Value = 145800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("D66").Value = 831500
$ws.Range("E66").Value = 953800
$ws.Range("F66").Value = 930200
$ws.Range("G66").Value = 854600
$ws.Range("H66").Value = 777700
$ws.Range("I66").Value = 725600
$ws.Range("J66").Value = 735000
$ws.Range("K66").Value = 854600
$ws.Range("L66").Value = 905200
$ws.Range("M66").Value = 840400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("D72").Value = -1318800
$ws.Range("E72").Value = -1502300
$ws.Range("F72").Value = -1560800
$ws.Range("G72").Value = -1564600
$ws.Range("H72").Value = -1581200
$ws.Range("I72").Value = -1587200
$ws.Range("J72").Value = -1593800
$ws.Range("K72").Value = -1739000
$ws.Range("L72").Value = -1792000
$ws.Range("M72").Value = -2432000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("D76").Value = 1489200
$ws.Range("E76").Value = 1311300
$ws.Range("F76").Value = 1271800
$ws.Range("G76").Value = 1241300
$ws.Range("H76").Value = 1191300
$ws.Range("I76").Value = 1173900
$ws.Range("J76").Value = 1202900
$ws.Range("K76").Value = 1173400
$ws.Range("L76").Value = 1122000
$ws.Range("M76").Value = 461800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43281
$ws.Range("G80").Value = 43190
$ws.Range("H80").Value = 43100
$ws.Range("I80").Value = 43008
$ws.Range("J80").Value = 42916
$ws.Range("K80").Value = 42825
$ws.Range("L80").Value = 42735
$ws.Range("M80").Value = 42643
$ws.Range("D81").Value = 185500
$ws.Range("E81").Value = 64700
$ws.Range("F81").Value = 9200
$ws.Range("G81").Value = 22100
$ws.Range("H81").Value = 11400
$ws.Range("I81").Value = 12000
$ws.Range("J81").Value = 96200
$ws.Range("K81").Value = 58600
$ws.Range("L81").Value = 645600
$ws.Range("M81").Value = -77400
$ws.Range("D83").Value = 63700
$ws.Range("E83").Value = 60700
$ws.Range("F83").Value = 54500
$ws.Range("G83").Value = 47700
$ws.Range("H83").Value = 48900
$ws.Range("I83").Value = 44500
$ws.Range("J83").Value = 48200
$ws.Range("K83").Value = 46500
$ws.Range("L83").Value = 49100
$ws.Range("M83").Value = 70400
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("D89").Value = 164900
$ws.Range("E89").Value = 160800
$ws.Range("F89").Value = 105500
$ws.Range("G89").Value = 118600
$ws.Range("H89").Value = 100700
$ws.Range("I89").Value = 85300
$ws.Range("J89").Value = 73200
$ws.Range("K89").Value = 98300
$ws.Range("L89").Value = 57400
$ws.Range("M89").Value = 81400
$ws.Range("D91").Value = -600
$ws.Range("E91").Value = 6300
$ws.Range("F91").Value = -1800
$ws.Range("G91").Value = -4600
$ws.Range("H91").Value = -2800
$ws.Range("I91").Value = -1700
$ws.Range("J91").Value = -4700
$ws.Range("K91").Value = -168800
$ws.Range("L91").Value = -432500
$ws.Range("M91").Value = -49700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("D94").Value = -109800
$ws.Range("E94").Value = -157400
$ws.Range("F94").Value = -112000
$ws.Range("G94").Value = -84500
$ws.Range("H94").Value = -381600
$ws.Range("I94").Value = -96600
$ws.Range("J94").Value = 239200
$ws.Range("K94").Value = -73000
$ws.Range("L94").Value = 174100
$ws.Range("M94").Value = -56500
$ws.Range("D96").Value = -5400
$ws.Range("E96").Value = -5500
$ws.Range("F96").Value = -5500
$ws.Range("G96").Value = -5400
$ws.Range("H96").Value = -5400
$ws.Range("I96").Value = -5400
$ws.Range("J96").Value = -5400
$ws.Range("K96").Value = -5600
$ws.Range("L96").Value = -5500
$ws.Range("M96").Value = -5500
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("D100").Value = -57600
$ws.Range("E100").Value = -8500
$ws.Range("F100").Value = -24900
$ws.Range("G100").Value = -4400
$ws.Range("H100").Value = -5400
$ws.Range("I100").Value = -5400
$ws.Range("J100").Value = -30100
$ws.Range("K100").Value = -20300
$ws.Range("L100").Value = 12300
$ws.Range("M100").Value = -5500
$ws.Range("D101").Value = 14200
$ws.Range("E101").Value = -4400
$ws.Range("F101").Value = 4600
$ws.Range("G101").Value = 7400
$ws.Range("H101").Value = -7300
$ws.Range("I101").Value = -10100
$ws.Range("J101").Value = 1200
$ws.Range("K101").Value = -2700
$ws.Range("L101").Value = 700
$ws.Range("M101").Value = 500
$ws.Range("D102").Value = 11700
$ws.Range("E102").Value = -9500
$ws.Range("F102").Value = -26800
$ws.Range("G102").Value = 37100
$ws.Range("H102").Value = -293600
$ws.Range("I102").Value = -26800
$ws.Range("J102").Value = 283500
$ws.Range("K102").Value = 2200
$ws.Range("L102").Value = 244500
$ws.Range("M102").Value = 19800
